$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price (D) and 1h volume change (E) columns with refreshed
# values from the symbol list feed. Values must remain plain text (as in
# the source data) rather than being auto-converted to numbers/percentages,
# so each value is written with a leading apostrophe and the cell style is
# reset to Normal afterwards (clears the quote-prefix flag Excel adds).

$ws.Range("D2").Value = "'301.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.57%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'0.95%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.090"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.77%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07850"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.94%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.336"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-6.37%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.814"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.48%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.840"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.08%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9167"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.78%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1757"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07544"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.64%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09228"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'15.63%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02994"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.76%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.64%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001517"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.47%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005798"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.42%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.471"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.68%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.250"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.38%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-1.14%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1328"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.40%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.050"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-11.83%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'11.60%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04622"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.79%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001249"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.27%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004462"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.40%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'5.70%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-1.53%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01762"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.90%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04744"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.38%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007348"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.55%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1358"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.28%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-2.45%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009780"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-8.59%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.50%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.16%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'20.08%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.7456"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-9.14%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.16%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.16%"
$ws.Range("E50").Style = "Normal"

Write-Host "Applied 66 cell updates"
